$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# Rename the first sheet
$ws1.Name = "火灾自动报警"

# Clear the "不限" (unlimited) values from the point-limit column on rows
# that no longer carry a cap
$ws1.Range("C4").ClearContents()
$ws1.Range("C6").ClearContents()
$ws1.Range("C8").ClearContents()
$ws1.Range("C9").ClearContents()
$ws1.Range("C10").ClearContents()
$ws1.Range("C11").ClearContents()

# Update the selection/active state on the "照明" sheet first ...
$ws2.Activate()
$ws2.Range("B2").Select()

# ... then switch to and select a cell on the renamed sheet so it ends up
# being the active tab when the workbook is saved
$ws1.Activate()
$ws1.Range("C35").Select()
